# Update InsideBet Data: Automatizado
# Applies updated stats for rows 9 (Estrela) and 18 (Tondela) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 - Estrela
$ws.Range("D9").Value = 47.8
$ws.Range("E9").Value = 23
$ws.Range("F9").Value = 253
$ws.Range("G9").Value = 2070
$ws.Range("H9").Value = 23
$ws.Range("O9").Value = 61
$ws.Range("P9").Value = 4
$ws.Range("Q9").Value = 1.13
$ws.Range("R9").Value = 0.78
$ws.Range("S9").Value = 1.91
$ws.Range("T9").Value = 0.96
$ws.Range("U9").Value = 1.74

# Row 18 - Tondela
$ws.Range("D18").Value = 41.6
$ws.Range("E18").Value = 23
$ws.Range("F18").Value = 253
$ws.Range("G18").Value = 2070
$ws.Range("H18").Value = 23
$ws.Range("I18").Value = 17
$ws.Range("K18").Value = 25
$ws.Range("L18").Value = 12
$ws.Range("M18").Value = 5
$ws.Range("N18").Value = 9
$ws.Range("O18").Value = 65
$ws.Range("Q18").Value = 0.74
$ws.Range("R18").Value = 0.35
$ws.Range("S18").Value = 1.09
$ws.Range("T18").Value = 0.52
$ws.Range("U18").Value = 0.87
